$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.698.26"
$ws.Range("E2").Value = "  -2.98%  "
$ws.Range("D3").Value = "1.852.82"
$ws.Range("E3").Value = "  -3.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -1.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.52"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4656"
$ws.Range("E7").Value = "  -3.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3916"
$ws.Range("E8").Value = "  -3.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.46"
$ws.Range("E9").Value = "  -2.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07927"
$ws.Range("E10").Value = "  -3.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9859"
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.37"
$ws.Range("E12").Value = "  -5.36%  "
$ws.Range("D13").Value = "1.901.56"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.850"
$ws.Range("E14").Value = "  -3.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.001"
$ws.Range("E15").Value = "  -3.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06839"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.76"
$ws.Range("E17").Value = "  -4.24%  "
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001006"
$ws.Range("E19").Value = "  -3.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.11"
$ws.Range("E20").Value = "  -2.86%  "
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "28.701.65"
$ws.Range("E22").Value = "  -2.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.394"
$ws.Range("E23").Value = "  -5.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.33"
$ws.Range("E24").Value = "  -4.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.131"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("D26").Value = "2.126.68"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.21"
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.44"
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.130"
$ws.Range("E29").Value = "  -4.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.025"
$ws.Range("E30").Value = "  -3.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.67"
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9763"
$ws.Range("E32").Value = "  -3.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09442"
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.371"
$ws.Range("E34").Value = "  -4.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.477"
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.354"
$ws.Range("E36").Value = "  -1.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06178"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.166"
$ws.Range("E39").Value = "  -1.71%  "
$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.001"
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5726"
$ws.Range("E41").Value = "  -3.89%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.617"
$ws.Range("E42").Value = "  -3.39%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.18"
$ws.Range("E43").Value = "  -5.46%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1799"
$ws.Range("E44").Value = "  -2.77%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.378"
$ws.Range("E45").Value = "  -3.03%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.249"
$ws.Range("E46").Value = "  -2.38%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5412"
$ws.Range("E47").Value = "  -2.79%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.79"
$ws.Range("E48").Value = "  -5.26%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07152"
$ws.Range("E49").Value = "  -4.50%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.914"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "114.15"
$ws.Range("E51").Value = "  -4.18%  "
